$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting old E (change_flag) onward to the right
$ws.Columns.Item(5).Insert()

# Header for the new column
$ws.Range("E1").Value = "tied_teams"

# Populate tied_teams values per row
$ws.Range("E2:E3").Value = '[''Romania'', ''Czech Republic'']'
$ws.Range("E4:E6").Value = '[]'
$ws.Range("E7:E9").Value = '[''Sweden'', ''Portugal'']'
$ws.Range("E10:E12").Value = '[''Northern Ireland'', ''Portugal'']'
$ws.Range("E13:E27").Value = '[]'
$ws.Range("E28:E32").Value = '[''Switzerland'', ''Hungary'']'
$ws.Range("E33:E35").Value = '[]'
$ws.Range("E36:E39").Value = '[''Denmark'', ''Slovakia'']'
$ws.Range("E40:E61").Value = '[]'
$ws.Range("E62:E63").Value = '[''Netherlands'', ''Slovakia'']'
$ws.Range("E64:E71").Value = '[]'
$ws.Range("E72:E73").Value = '[''Netherlands'', ''Slovakia'']'
$ws.Range("E74:E75").Value = '[]'
$ws.Range("E76:E80").Value = '[''Netherlands'', ''Georgia'']'
